$d = $word.ActiveDocument
$wOpen = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'

# Locate the two section headings the rest of the edit is anchored to.
$resistIndex = 0
$fusionIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $text = $d.Paragraphs($i).Range.Text
    if ($text -like "*Resist*ncia:*") { $resistIndex = $i }
    if ($text -like "*Modelos no Fusion:*") { $fusionIndex = $i }
}

# 1. Remove the existing "_GoBack" bookmark from the Space Shuttle paragraph
#    (it sat between the "...exaustores" run and the ", forças de arrasto..." run).
$d.Bookmarks("_GoBack").Delete()

# 2. Re-add the "_GoBack" bookmark, now spanning the second of the two blank
#    paragraphs that follow "Termômetro de Resistência:".
$goBackPara = $d.Paragraphs($resistIndex + 2)
$goBackXml = $wOpen + '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$goBackPara.Range.InsertXML($goBackXml)

# 3. Fill in the first blank paragraph after "Modelos no Fusion:" with the
#    new descriptive text (including the proofing marks Word leaves around
#    the words it doesn't recognise).
$fusionPara1 = $d.Paragraphs($fusionIndex + 1)
$fusionXml1 = $wOpen + '<w:r><w:t xml:space="preserve">A proposta era criar 14 modelos diferentes no </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>program</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> do </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>AutoDesk</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> Fusion. Além de desenhá-los</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">com representação e escalas </w:t></w:r></w:p>'
$fusionPara1.Range.InsertXML($fusionXml1)

# 4. Turn the second blank paragraph after "Modelos no Fusion:" into an empty
#    paragraph whose mark is tagged as English (en-US) run formatting.
$fusionPara2 = $d.Paragraphs($fusionIndex + 2)
$fusionXml2 = $wOpen + '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>'
$fusionPara2.Range.InsertXML($fusionXml2)
